# Add a new day row (19) to the workload tracker, describing the work done on
# decimal-to-fraction conversion and finding the angle between clock hands.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a copy of the previous day's row (18) so the new row inherits the
# exact same cell formatting/styles (date format, borders, alignment, etc.).
$ws.Range("A18:E18").Copy($ws.Range("A19:E19"))

# Day for the new row: 2013-02-21 (serial date 41326), one day after row 18.
$ws.Range("A19").Value = 41326

# Task durations for the new day.
$ws.Range("B19").Value = "0H"
$ws.Range("C19").Value = "1H"
$ws.Range("D19").Value = "0.5H"

# Description of the math exercise worked on.
$ws.Range("E19").Value = "Decimal to Fraction"

# Move the active selection to the newly added row, mirroring what Excel does
# after entering data in a new row.
[void]$ws.Range("A19").Select()
